$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Append new row 11 for "Gobernatura 23" first (so its strings register earlier
# in the shared string table), then insert the "Senado 18" row above it.
$ws.Cells.Item(11, 1).Value = "Gobernatura 23"
$ws.Cells.Item(11, 2).Value = "gb_23"
$ws.Cells.Item(11, 3).Value = "#283618"

# Insert a new row at position 8 for "Senado 18" (shifts existing rows 8-11 down to 9-12)
$ws.Rows.Item(8).Insert()

$ws.Cells.Item(8, 1).Value = "Senado 18"
$ws.Cells.Item(8, 2).Value = "sen_18"
$ws.Cells.Item(8, 3).Value = "348cae4"

# Set column A width (closest achievable value to the target 30.1640625)
$ws.Columns.Item(1).ColumnWidth = 29.25

# Update selection to C8
$ws.Range("C8").Select()
